$d = $word.ActiveDocument

# Locate the paragraph that ends the "row/column" answer text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*A column is a vertical series of cells in a chart, table, or spreadsheet.*") {
        $target = $p
        break
    }
}

# Insert a brand-new empty paragraph right after it.
$newPara = $target.Range.InsertParagraphAfter()

# Re-resolve the freshly-inserted paragraph and configure its formatting
# to match the diff: BodyText style, single-spaced, hanging indent at
# 720 twips (no first-line indent), Times New Roman 14pt (sz 28 half-points).
$newRange = $target.Next()
$newRange.Style = "Body Text"
$newRange.Format.SpaceAfter = 0
$newRange.Format.LineSpacingRule = 0
$newRange.Format.LineSpacing = 12
$newRange.Format.LeftIndent = 36
$newRange.Format.FirstLineIndent = 0
$newRange.Range.Font.Name = "Times New Roman"
$newRange.Range.Font.Size = 14
